$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PMOIY")

# Row 9 - Cost of Revenue
$ws.Range("D9").Value = 231600
$ws.Range("E9").Value = 236600
$ws.Range("F9").Value = 218800
$ws.Range("G9").Value = 405000
$ws.Range("H9").Value = 348000
$ws.Range("I9").Value = 362200
$ws.Range("J9").Value = 298800

# Row 10 - Gross Profit
$ws.Range("D10").Value = 393400
$ws.Range("E10").Value = 260400
$ws.Range("F10").Value = 327300
$ws.Range("G10").Value = 164900
$ws.Range("H10").Value = 19100
$ws.Range("I10").Value = 128000
$ws.Range("J10").Value = 278200

# Row 91 - Capital Expenditures
$ws.Range("G91").Value = 309500
$ws.Range("H91").Value = -318300
